$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '71.103.47'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.851.67'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '695.64'
$ws.Range('E5').Value = '  +3.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.24'
$ws.Range('E6').Value = '  +2.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.851.28'
$ws.Range('E7').Value = '  +1.01%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +0.34%  '
$ws.Range('E10').Value = '  +1.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.31'
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.462'
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000258'
$ws.Range('E13').Value = '  +6.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.57'
$ws.Range('E14').Value = '  +1.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.506.75'
$ws.Range('E15').Value = '  +1.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.856.27'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '71.175.44'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.78'
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('E19').Value = '  +0.91%  '
$ws.Range('E20').Value = '  +0.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.13'
$ws.Range('E21').Value = '  -5.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '495.11'
$ws.Range('E22').Value = '  +3.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.724'
$ws.Range('E23').Value = '  +1.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.92'
$ws.Range('E24').Value = '  +1.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000147'
$ws.Range('E25').Value = '  +4.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.35'
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('E27').Value = '  +2.87%  '
$ws.Range('E28').Value = '  +2.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.009.54'
$ws.Range('E29').Value = '  +1.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.18'
$ws.Range('E30').Value = '  +10.83%  '
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('E32').Value = '  +3.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.30'
$ws.Range('E33').Value = '  +0.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.74'
$ws.Range('E34').Value = '  +0.75%  '
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.31'
$ws.Range('E36').Value = '  +2.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.808.22'
$ws.Range('E37').Value = '  +1.00%  '
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('E39').Value = '  +3.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.37'
$ws.Range('E40').Value = '  +12.14%  '
$ws.Range('E41').Value = '  +1.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.06'
$ws.Range('E42').Value = '  +2.06%  '
$ws.Range('E43').Value = '  +6.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '163.84'
$ws.Range('E46').Value = '  +2.55%  '
$ws.Range('E47').Value = '  +5.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '48.67'
$ws.Range('E48').Value = '  +1.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '44.34'
$ws.Range('E49').Value = '  -3.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '420.42'
$ws.Range('E50').Value = '  +6.70%  '
$ws.Range('E51').Value = '  +0.91%  '
